# Update attendance/view numbers (column F) in the "展览" and "全部类型"
# worksheets to reflect the latest scrape, per commit:
# "Update gh-pages to output generated at 456a3b4"

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value  = 1756
$ws1.Range("F11").Value = 1952
$ws1.Range("F13").Value = 607
$ws1.Range("F14").Value = 449
$ws1.Range("F15").Value = 12
$ws1.Range("F18").Value = 7
$ws1.Range("F23").Value = 1050
$ws1.Range("F24").Value = 4
$ws1.Range("F25").Value = 322
$ws1.Range("F27").Value = 260
$ws1.Range("F28").Value = 301

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value  = 1756
$ws4.Range("F12").Value = 1952
$ws4.Range("F14").Value = 607
$ws4.Range("F15").Value = 449
$ws4.Range("F16").Value = 12
$ws4.Range("F19").Value = 7
$ws4.Range("F24").Value = 1050
$ws4.Range("F25").Value = 4
$ws4.Range("F26").Value = 322
$ws4.Range("F28").Value = 260
$ws4.Range("F29").Value = 301
